# use-case-pesquisar-expressao.xlsx edit script
# Implements: rename "Pesquisar palavra/expressao" -> "Pesquisar expressao" use case,
# rewrite the normal-scenario steps, shrink the pre-condition text/row height,
# and add a new "Excecao 1" box (rows 11-12) describing the invalid-expression path.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header block (rows 2-5): update the use case title and pre-condition text
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Pesquisar expressão"
$ws.Range("C4").Value = "Registado/logado no sistema, `nter escolhido e iniciado uma receita"

# row 4 used to be tall enough for a 3-line condition; new text is shorter
$ws.Rows.Item(4).RowHeight = 40

# ---------------------------------------------------------------------------
# 2. Normal scenario steps (rows 7-10)
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "1. Indicar expressão para explicação"

$ws.Range("D8").Value = "2. Validar expressão"

$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = "3. Apresenta explicação da expressão"

$ws.Range("D10").Value = "4. Termina processo"
$ws.Rows.Item(8).RowHeight = 20

# ---------------------------------------------------------------------------
# 3. New "Exceção 1" box (rows 11-12) under the normal scenario table
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(12).Insert()

# Reuse formatting from the "Cenário Normal" header row (grey label cell +
# plain input/output cells) as a starting point, then fix up borders below.
$ws.Range("B6:D6").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122)
$ws.Range("B6:D6").Copy()
$ws.Range("B12:D12").PasteSpecial(-4122)

$ws.Range("B11:B12").Merge()

$ws.Range("B11").Value = "Exceção 1`n(passo 2)`n[expressão inválida]`n"
$ws.Range("D11").Value = "2.1. Informa que a expressão é inválida"
$ws.Range("D12").Value = "2.2 Termina processo "

$ws.Range("B11:B12").Font.Size = 14
$ws.Range("C11:D12").Font.Size = 14
$ws.Range("B11:B12").HorizontalAlignment = -4108
$ws.Range("B11:B12").VerticalAlignment = -4160
$ws.Range("B11").WrapText = $true
$ws.Range("D12").VerticalAlignment = -4160

$ws.Rows.Item(11).RowHeight = 20
$ws.Rows.Item(12).RowHeight = 60

# Borders around B11:D12 (outer medium box, medium divider under row11 in
# column B only continues as part of the merge, thin divider between C/D
# and the label column, medium divider between row 11 and row 12 in C:D)
$outer = $ws.Range("B11:D12")
$outer.Borders.Item(7).LineStyle = 1
$outer.Borders.Item(7).Weight = -4138
$outer.Borders.Item(10).LineStyle = 1
$outer.Borders.Item(10).Weight = -4138
$outer.Borders.Item(8).LineStyle = 1
$outer.Borders.Item(8).Weight = -4138
$outer.Borders.Item(9).LineStyle = 1
$outer.Borders.Item(9).Weight = -4138

$bc = $ws.Range("B11:C12")
$bc.Borders.Item(11).LineStyle = 1
$bc.Borders.Item(11).Weight = 2

$cd = $ws.Range("C11:C12")
$cd.Borders.Item(12).LineStyle = 1
$cd.Borders.Item(12).Weight = -4138

Write-Output "done"
